# RAP3.xlsx multistage-build edit
# - removes the incomplete "Esther" account row (old row 21)
# - fixes Jan's password (F column) to match the userid
# - adds a "label"/"PF_Label" column to the Roles lookup table
# - appends 6 new system roles (SystemAdmin, Administrator, Anonymous,
#   ExecEngine, Janitor, SYSTEM) with a distinct (Menlo) font

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the old row 21 (incomplete "Esther" account entry). This
#    shifts every row below it up by one and Excel automatically fixes
#    up the relative formulas (CONCATENATE("Acc_",$B7) etc.)
$ws.Rows("21").Delete()

# 2. Fix Jan's password in column F (row 21 after the shift) so it
#    matches his userid in B21 (123456, not 12345)
$ws.Range("F21").Value = 123456

# 3. Add the "label" column (C) to the Roles lookup table
$ws.Range("C23").Value = "label"
$ws.Range("C24").Value = "PF_Label"
$ws.Range("C25").Value = "Tutor"
$ws.Range("C26").Value = "Student"
$ws.Range("C27").Value = "GradStudent"
$ws.Range("C28").Value = "AccountManager"

# match the formatting of the existing B column cells on those rows
$ws.Range("B23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("B24").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("B25").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("B26").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("B27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("B28").Copy()
$ws.Range("C28").PasteSpecial(-4122)

# 4. Append the 6 new system roles below the existing ones
$ws.Range("A29").Value = "SystemAdmin"
$ws.Range("C29").Value = "SystemAdmin"

$ws.Range("A30").Value = "Administrator"
$ws.Range("B30").Value = "Administrator"

$ws.Range("A31").Value = "Anonymous"
$ws.Range("B31").Value = "Anonymous"

$ws.Range("A32").Value = "ExecEngine"
$ws.Range("B32").Value = "ExecEngine"

$ws.Range("A33").Value = "Janitor"
$ws.Range("B33").Value = "Janitor"

$ws.Range("A34").Value = "SYSTEM"
$ws.Range("B34").Value = "SYSTEM"

# give the new rows the distinct Menlo font used for these system roles
$newRolesFont = $ws.Range("A29:C34")
$newRolesFont.Font.Name = "Menlo"
$newRolesFont.Font.Color = 0

# 5. Update the selection to match the saved workbook state
$ws.Range("K21").Select()
